$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1242
$ws.Cells.Item(32, 9).Value = 865.5833
$ws.Cells.Item(32, 10).Value = 3500.5
$ws.Cells.Item(32, 11).Value = 865.5833
$ws.Cells.Item(32, 12).Value = 3500.5
$ws.Cells.Item(32, 13).Value = -539.5833
$ws.Cells.Item(32, 14).Value = -4152.5
$ws.Cells.Item(33, 8).Value = 141
$ws.Cells.Item(33, 9).Value = 145.66667
$ws.Cells.Item(33, 11).Value = 145.66667
$ws.Cells.Item(33, 13).Value = 83.33332999999999
$ws.Cells.Item(51, 8).Value = 4999.5
$ws.Cells.Item(51, 9).Value = 4999.5
$ws.Cells.Item(51, 11).Value = 4999.5
$ws.Cells.Item(51, 13).Value = -4515.5
$ws.Cells.Item(70, 8).Value = 148624.2
$ws.Cells.Item(70, 10).Value = 246373.67
$ws.Cells.Item(70, 12).Value = 739121.01
$ws.Cells.Item(70, 14).Value = -739661.01
$ws.Cells.Item(73, 8).Value = 148624.2
$ws.Cells.Item(73, 10).Value = 246373.67
$ws.Cells.Item(73, 12).Value = 739121.01
$ws.Cells.Item(73, 14).Value = -740993.01
$ws.Cells.Item(132, 8).Value = 888.13336
$ws.Cells.Item(132, 9).Value = 915.4643
$ws.Cells.Item(132, 11).Value = 2746.3929
$ws.Cells.Item(132, 13).Value = -216.3928999999998
$ws.Cells.Item(138, 8).Value = 4838.523
$ws.Cells.Item(138, 10).Value = 5608.0264
$ws.Cells.Item(138, 12).Value = 16824.0792
$ws.Cells.Item(138, 14).Value = -27104.0792
$ws.Cells.Item(141, 8).Value = 1887.1904
$ws.Cells.Item(141, 9).Value = 1814.9
$ws.Cells.Item(141, 11).Value = 5444.700000000001
$ws.Cells.Item(141, 13).Value = -264.7000000000007

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5258.0977
$ws.Cells.Item(32, 9).Value = 3123.838
$ws.Cells.Item(32, 11).Value = 3123.838
$ws.Cells.Item(32, 13).Value = -2836.838
$ws.Cells.Item(39, 8).Value = 504000
$ws.Cells.Item(39, 9).Value = 504000
$ws.Cells.Item(39, 11).Value = 504000
$ws.Cells.Item(39, 13).Value = -503480
$ws.Cells.Item(45, 8).Value = 799.5
$ws.Cells.Item(45, 9).Value = 799.5
$ws.Cells.Item(45, 11).Value = 799.5
$ws.Cells.Item(45, 13).Value = -422.5
$ws.Cells.Item(124, 8).Value = 42357
$ws.Cells.Item(124, 10).Value = 42357
$ws.Cells.Item(124, 12).Value = 42357
$ws.Cells.Item(124, 14).Value = -52177
$ws.Cells.Item(132, 8).Value = 788.43475
$ws.Cells.Item(132, 9).Value = 693.1429000000001
$ws.Cells.Item(132, 11).Value = 2079.4287
$ws.Cells.Item(132, 13).Value = 450.5712999999996

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3562.4167
$ws.Cells.Item(105, 9).Value = 3878.4285
$ws.Cells.Item(105, 11).Value = 3878.4285
$ws.Cells.Item(105, 13).Value = -2131.4285

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4571.1577
$ws.Cells.Item(31, 9).Value = 3979.6428
$ws.Cells.Item(31, 10).Value = 6227.4
$ws.Cells.Item(31, 11).Value = 3979.6428
$ws.Cells.Item(31, 12).Value = 6227.4
$ws.Cells.Item(31, 13).Value = -3684.6428
$ws.Cells.Item(31, 14).Value = -6817.4
$ws.Cells.Item(34, 8).Value = 4571.1577
$ws.Cells.Item(34, 9).Value = 3979.6428
$ws.Cells.Item(34, 10).Value = 6227.4
$ws.Cells.Item(34, 11).Value = 3979.6428
$ws.Cells.Item(34, 12).Value = 6227.4
$ws.Cells.Item(34, 13).Value = -3777.6428
$ws.Cells.Item(34, 14).Value = -6631.4
$ws.Cells.Item(39, 8).Value = 2300
$ws.Cells.Item(39, 9).Value = 2300
$ws.Cells.Item(39, 11).Value = 2300
$ws.Cells.Item(39, 13).Value = -1909
$ws.Cells.Item(49, 8).Value = 2300
$ws.Cells.Item(49, 9).Value = 2300
$ws.Cells.Item(49, 11).Value = 2300
$ws.Cells.Item(49, 13).Value = -2118
$ws.Cells.Item(107, 8).Value = 1858.5
$ws.Cells.Item(107, 10).Value = 2950
$ws.Cells.Item(107, 12).Value = 2950
$ws.Cells.Item(107, 14).Value = -6790
$ws.Cells.Item(109, 8).Value = 42564.75
$ws.Cells.Item(109, 10).Value = 40000
$ws.Cells.Item(109, 12).Value = 40000
$ws.Cells.Item(109, 14).Value = -42080
$ws.Cells.Item(132, 8).Value = 1582.35
$ws.Cells.Item(132, 9).Value = 1591.5555
$ws.Cells.Item(132, 10).Value = 1499.5
$ws.Cells.Item(132, 11).Value = 4774.666499999999
$ws.Cells.Item(132, 12).Value = 4498.5
$ws.Cells.Item(132, 13).Value = -2244.666499999999
$ws.Cells.Item(132, 14).Value = -9558.5
$ws.Cells.Item(134, 8).Value = 2739.652
$ws.Cells.Item(134, 9).Value = 2674.2778
$ws.Cells.Item(134, 11).Value = 8022.8334
$ws.Cells.Item(134, 13).Value = -5487.8334

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(120, 8).Value = 13366.667
$ws.Cells.Item(120, 9).Value = 9100
$ws.Cells.Item(120, 10).Value = 15500
$ws.Cells.Item(120, 11).Value = 27300
$ws.Cells.Item(120, 12).Value = 46500
$ws.Cells.Item(120, 13).Value = -22462
$ws.Cells.Item(120, 14).Value = -56176
$ws.Cells.Item(131, 8).Value = 1310.4067
$ws.Cells.Item(131, 9).Value = 929
$ws.Cells.Item(131, 10).Value = 1330.8392
$ws.Cells.Item(131, 11).Value = 2787
$ws.Cells.Item(131, 12).Value = 3992.5176
$ws.Cells.Item(131, 13).Value = 2253
$ws.Cells.Item(131, 14).Value = -14072.5176

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 553181.25
$ws.Cells.Item(24, 9).Value = 3000000
$ws.Cells.Item(24, 11).Value = 3000000
$ws.Cells.Item(24, 13).Value = -2999827
$ws.Cells.Item(46, 8).Value = 17691.846
$ws.Cells.Item(46, 9).Value = 9998
$ws.Cells.Item(46, 10).Value = 20000
$ws.Cells.Item(46, 11).Value = 9998
$ws.Cells.Item(46, 12).Value = 20000
$ws.Cells.Item(46, 13).Value = -9842
$ws.Cells.Item(46, 14).Value = -20312
$ws.Cells.Item(80, 8).Value = 10157.467
$ws.Cells.Item(80, 9).Value = 4126.75
$ws.Cells.Item(80, 11).Value = 4126.75
$ws.Cells.Item(80, 13).Value = -3128.75
$ws.Cells.Item(83, 8).Value = 10157.467
$ws.Cells.Item(83, 9).Value = 4126.75
$ws.Cells.Item(83, 11).Value = 20633.75
$ws.Cells.Item(83, 13).Value = -15641.75
$ws.Cells.Item(107, 8).Value = 1000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(122, 8).Value = 103253.3
$ws.Cells.Item(122, 9).Value = 2937.8572
$ws.Cells.Item(122, 10).Value = 337322.66
$ws.Cells.Item(122, 11).Value = 8813.571599999999
$ws.Cells.Item(122, 12).Value = 1011967.98
$ws.Cells.Item(122, 13).Value = -6363.571599999999
$ws.Cells.Item(122, 14).Value = -1016867.98
$ws.Cells.Item(123, 8).Value = 51250.375
$ws.Cells.Item(123, 10).Value = 51250.375
$ws.Cells.Item(123, 12).Value = 51250.375
$ws.Cells.Item(123, 14).Value = -56150.375
$ws.Cells.Item(132, 8).Value = 2833.75
$ws.Cells.Item(132, 9).Value = 2612.0557
$ws.Cells.Item(132, 11).Value = 7836.1671
$ws.Cells.Item(132, 13).Value = -5306.1671

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3600
$ws.Cells.Item(136, 9).Value = 3799.5
$ws.Cells.Item(136, 11).Value = 11398.5
$ws.Cells.Item(136, 13).Value = -8848.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 5167.3335
$ws.Cells.Item(122, 9).Value = 5999.75
$ws.Cells.Item(122, 11).Value = 17999.25
$ws.Cells.Item(122, 13).Value = -15549.25
$ws.Cells.Item(132, 8).Value = 3668.1875
$ws.Cells.Item(132, 9).Value = 715.1539
$ws.Cells.Item(132, 10).Value = 16464.666
$ws.Cells.Item(132, 11).Value = 2145.4617
$ws.Cells.Item(132, 12).Value = 49393.99800000001
$ws.Cells.Item(132, 13).Value = 384.5383000000002
$ws.Cells.Item(132, 14).Value = -54453.99800000001

Write-Output "done"